# Apply weekly update to Haba (Hortaliza) data sheet.
# The edit is a permutation of the data rows (2-14): the values in columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are rearranged across
# the existing rows. All other columns stay the same for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (after the edit), keyed by row number.
$data = @{
    2  = @{ D = 44476; J = 900;  K = 700;  L = 800;  M = 750;  P = 750 }
    3  = @{ D = 44484; J = 900;  K = 750;  L = 800;  M = 775;  P = 775 }
    4  = @{ D = 44243; J = 1200; K = 1200; L = 1300; M = 1250; P = 1250 }
    5  = @{ D = 44455; J = 1100; K = 900;  L = 1000; M = 950;  P = 950 }
    6  = @{ D = 44449; J = 1300; K = 900;  L = 950;  M = 925;  P = 925 }
    7  = @{ D = 44442; J = 1250; K = 850;  L = 900;  M = 875;  P = 875 }
    8  = @{ D = 44175; J = 1600; K = 1000; L = 1200; M = 1100; P = 1100 }
    9  = @{ D = 44291; J = 1000; K = 1000; L = 1200; M = 1100; P = 1100 }
    10 = @{ D = 44453; J = 1000; K = 800;  L = 900;  M = 850;  P = 850 }
    11 = @{ D = 44229; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    12 = @{ D = 44407; J = 1000; K = 1200; L = 1300; M = 1250; P = 1250 }
    13 = @{ D = 44284; J = 1500; K = 800;  L = 850;  M = 825;  P = 825 }
    14 = @{ D = 44341; J = 1300; K = 900;  L = 1000; M = 950;  P = 950 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # Column D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # Column J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # Column K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # Column L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # Column M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals.P   # Column P - Precio $/Kg
}
